# Apply the updated cryptos list values captured on Thu Mar  9 23:56:42 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without leaving a
# lingering "stored as text" number-format style behind (mirrors the
# original inline-string cells, which carry no cell style).
function Set-TextValue($address, $text) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '20.387.28'
Set-TextValue 'E2' '  -6.42%  '
# Row 3
Set-TextValue 'D3' '1.440.62'
Set-TextValue 'E3' '  -6.41%  '
# Row 4
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  -0.34%  '
# Row 5
Set-TextValue 'E5' '  -0.25%  '
# Row 6
Set-TextValue 'D6' '277.96'
Set-TextValue 'E6' '  -3.43%  '
# Row 7
Set-TextValue 'D7' '0.3728'
Set-TextValue 'E7' '  -4.42%  '
# Row 8
Set-TextValue 'D8' '0.3101'
Set-TextValue 'E8' '  -2.63%  '
# Row 9
Set-TextValue 'D9' '40.69'
Set-TextValue 'E9' '  -5.36%  '
# Row 10
Set-TextValue 'D10' '1.017'
Set-TextValue 'E10' '  -3.60%  '
# Row 11
Set-TextValue 'D11' '0.06607'
Set-TextValue 'E11' '  -7.36%  '
# Row 12
Set-TextValue 'D12' '1.002'
Set-TextValue 'E12' '  -0.34%  '
# Row 13
Set-TextValue 'D13' '5.389'
Set-TextValue 'E13' '  -4.08%  '
# Row 14
Set-TextValue 'D14' '17.36'
Set-TextValue 'E14' '  -6.09%  '
# Row 15
Set-TextValue 'D15' '6.169'
Set-TextValue 'E15' '  -6.65%  '
# Row 16
Set-TextValue 'D16' '1.439.06'
Set-TextValue 'E16' '  -7.09%  '
# Row 17
Set-TextValue 'D17' '0.00001010'
Set-TextValue 'E17' '  -8.06%  '
# Row 18
Set-TextValue 'D18' '76.81'
Set-TextValue 'E18' '  -7.58%  '
# Row 19
Set-TextValue 'D19' '0.05858'
Set-TextValue 'E19' '  -10.62%  '
# Row 20
Set-TextValue 'E20' '  -0.22%  '
# Row 21
Set-TextValue 'D21' '5.748'
Set-TextValue 'E21' '  -6.14%  '
# Row 22
Set-TextValue 'D22' '14.44'
Set-TextValue 'E22' '  -4.97%  '
# Row 23
Set-TextValue 'D23' '11.05'
Set-TextValue 'E23' '  +0.45%  '
# Row 24
Set-TextValue 'D24' '2.323'
Set-TextValue 'E24' '  -3.11%  '
# Row 25
Set-TextValue 'D25' '20.375.79'
Set-TextValue 'E25' '  -6.55%  '
# Row 26
Set-TextValue 'D26' '2.280'
Set-TextValue 'E26' '  -3.72%  '
# Row 27
Set-TextValue 'D27' '142.24'
Set-TextValue 'E27' '  -1.55%  '
# Row 28
Set-TextValue 'D28' '17.09'
Set-TextValue 'E28' '  -6.83%  '
# Row 29
Set-TextValue 'D29' '1.601.18'
Set-TextValue 'E29' '  -7.07%  '
# Row 30
Set-TextValue 'D30' '110.44'
Set-TextValue 'E30' '  -5.58%  '
# Row 31
Set-TextValue 'D31' '3.960'
Set-TextValue 'E31' '  -18.38%  '
# Row 32
Set-TextValue 'D32' '0.9311'
Set-TextValue 'E32' '  -3.29%  '
# Row 33
Set-TextValue 'D33' '5.499'
Set-TextValue 'E33' '  -5.43%  '
# Row 34
Set-TextValue 'D34' '0.07734'
Set-TextValue 'E34' '  -5.70%  '
# Row 35
Set-TextValue 'D35' '8.397'
Set-TextValue 'E35' '  -6.77%  '
# Row 36
Set-TextValue 'D36' '11.11'
Set-TextValue 'E36' '  +5.11%  '
# Row 37
Set-TextValue 'D37' '0.05756'
Set-TextValue 'E37' '  -5.30%  '
# Row 38
Set-TextValue 'E38' '  -0.24%  '
# Row 39
Set-TextValue 'D39' '4.759'
Set-TextValue 'E39' '  -6.28%  '
# Row 40
Set-TextValue 'B40' 'Algorand'
Set-TextValue 'C40' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.1936'
Set-TextValue 'E40' '  -4.56%  '
# Row 41
Set-TextValue 'B41' 'TrustWalletToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.133'
Set-TextValue 'E41' '  -5.02%  '
# Row 42
Set-TextValue 'D42' '0.02032'
Set-TextValue 'E42' '  -8.43%  '
# Row 43
Set-TextValue 'D43' '1.344'
Set-TextValue 'E43' '  -10.25%  '
# Row 44
Set-TextValue 'D44' '3.595'
Set-TextValue 'E44' '  -3.66%  '
# Row 45
Set-TextValue 'D45' '0.5355'
Set-TextValue 'E45' '  -6.16%  '
# Row 46
Set-TextValue 'D46' '12.13'
Set-TextValue 'E46' '  -5.43%  '
# Row 47
Set-TextValue 'D47' '0.5189'
Set-TextValue 'E47' '  -5.62%  '
# Row 48
Set-TextValue 'D48' '112.21'
Set-TextValue 'E48' '  -3.42%  '
# Row 49
Set-TextValue 'D49' '1.791'
Set-TextValue 'E49' '  -2.98%  '
# Row 50
Set-TextValue 'D50' '1.059'
Set-TextValue 'E50' '  -5.75%  '
# Row 51
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.06277'
Set-TextValue 'E51' '  -6.92%  '

Write-Host "Updated cryptos list values"
